# Fix mistake where combustion and biomass dataframes were mixed up
# (df_comb / df_biomass filtering bug) for several sheets.

$wb = $excel.ActiveWorkbook

# --- Strausberg: single value fix ---
$wsStrausberg = $wb.Worksheets.Item("Strausberg")
$wsStrausberg.Range("C3").Value = 3199

# --- Bocholt: row 2 values fix ---
$wsBocholt = $wb.Worksheets.Item("Bocholt")
$wsBocholt.Range("C2").Value = 6475
$wsBocholt.Range("D2").Value = 12655
$wsBocholt.Range("E2").Value = 11

# --- Kassel: rows 2-3 replaced with correct combustion-engine data,
#     and the (duplicate/incorrect) row 4 removed entirely ---
$wsKassel = $wb.Worksheets.Item("Kassel")

$wsKassel.Range("A2").Value = "Verbrennungsmotor"
$wsKassel.Range("B2").Value = "Biogas"
$wsKassel.Range("C2").Value = 806
$wsKassel.Range("D2").Value = 3108
$wsKassel.Range("E2").Value = 3

$wsKassel.Range("A3").Value = "Verbrennungsmotor"
$wsKassel.Range("B3").Value = "Biomethan (Bioerdgas)"
$wsKassel.Range("C3").Value = 3745
$wsKassel.Range("D3").Value = 804
$wsKassel.Range("E3").Value = 2

# Delete row 4 (its content is no longer valid; table shrinks to A1:E3)
$wsKassel.Rows.Item(4).Delete()
